$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 5-10: fill in "yes" in column E to match D/F/G, and set I/K numeric scores
$ws.Range("E5").Value = "yes"
$ws.Range("E6").Value = "yes"
$ws.Range("E7").Value = "yes"
$ws.Range("E8").Value = "yes"
$ws.Range("E9").Value = "yes"
$ws.Range("E10").Value = "yes"

$ws.Range("I5").Value = 4
$ws.Range("I6").Value = 4
$ws.Range("I7").Value = 4
$ws.Range("I8").Value = 4
$ws.Range("I9").Value = 4
$ws.Range("I10").Value = 4

$ws.Range("K5").Value = 4
$ws.Range("K6").Value = 4
$ws.Range("K7").Value = 4
$ws.Range("K8").Value = 4
$ws.Range("K9").Value = 3
$ws.Range("K10").Value = 3

# Rows 24-28: fill in "yes" in columns F and G, and set K score
$ws.Range("F24").Value = "yes"
$ws.Range("G24").Value = "yes"
$ws.Range("F25").Value = "yes"
$ws.Range("G25").Value = "yes"
$ws.Range("F26").Value = "yes"
$ws.Range("G26").Value = "yes"
$ws.Range("F27").Value = "yes"
$ws.Range("G27").Value = "yes"
$ws.Range("F28").Value = "yes"
$ws.Range("G28").Value = "yes"

$ws.Range("K24").Value = 2.5
$ws.Range("K25").Value = 2.5
$ws.Range("K26").Value = 2.5
$ws.Range("K27").Value = 2.5
$ws.Range("K28").Value = 2.5

# Update the view's selection / scroll position to match the saved workbook state
$ws.Application.ActiveWindow.ScrollColumn = 14
$ws.Range("H27").Select()
